# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G (header "K") values for rows 2-50 are recalculated/rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0, 0, 1, 1, 2, 0, 2, 0, 2, 1, 2, 3, 0, 4, 2, 1, 0, 3, 2, 0, 2, 0, 3, 0, 1, 0, 2, 1, 2, 1, 2, 0, 0, 1, 1, 0, 1, 0, 2, 3, 0, 2, 1, 3, 1, 1, 2, 0, 1)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
